$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New data rows appended to the daily user-impact status table.
$newRows = @(
    @(46049, 5605, 4474, 4035, 282, 78, 66, 12, 1),
    @(46050, 5604, 4478, 4042, 289, 73, 60, 12, 2)
)

$startRow = 78
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $values = $newRows[$i]

    # Column A holds a date (serial number), formatted like the rows above it.
    $ws.Cells.Item($r, 1).Value = $values[0]
    $ws.Cells.Item($r, 1).NumberFormat = $ws.Cells.Item($r - 1, 1).NumberFormat

    for ($c = 2; $c -le 9; $c++) {
        $ws.Cells.Item($r, $c).Value = $values[$c - 1]
    }
}

# Match the saved selection state from the edit.
$ws.Range("A79:I79").Select()
